$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("TestCredentials")
$ws1.Range("C3").Value = "authorprod@knowledgehut.com"
$ws1.Range("D3").Value = "Password@123"
$ws1.Activate()
$ws1.Range("C3").Select()
